$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.978990912437439
$ws.Range("B1").Value = 2.818645000457764
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.299554824829102
$ws.Range("E1").Value = 2.2932288646698
